# "add part buy and sale product"
# Rework Sheet1 from an employee table (ID/Name/Position/Address/Email/Phone/Parts/Manager)
# into a region/part table (ID/Name/Note) with new KV/SP rows.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New header for column C ("Position" -> "Note"); A1/B1 (ID/Name) stay as-is.
$ws1.Range("C1").Value = "Note"

# Row 2: Miền Nam region
$ws1.Range("A2").Value = "KV000001"
$ws1.Range("B2").Value = "Miền Nam"

# Row 3: second region row, also carries a note in column C
$ws1.Range("A3").Value = "KV000002"
$ws1.Range("B3").Value = "hieu"
$ws1.Range("C3").Value = "hieu"

# Row 4: Miền Trung region
$ws1.Range("A4").Value = "KV000003"
$ws1.Range("B4").Value = "Miền Trung"

# Row 5: Đông Nam Bộ region
$ws1.Range("A5").Value = "KV000004"
$ws1.Range("B5").Value = "Đông Nam Bộ"

# Row 6: new Part/product entry
$ws1.Range("A6").Value = "SP00005"
$ws1.Range("B6").Value = "Trung Du"

# Remove the old trailing employee columns (Address/Email/Phone/Parts/Manager)
# so the used range shrinks back down to A:C.
$ws1.Range("D1:H1").EntireColumn.Delete()

# Re-apply (approximate) autofit-style column widths for the new A:C layout.
$ws1.Columns.Item(1).ColumnWidth = 9.67
$ws1.Columns.Item(2).ColumnWidth = 13.17
$ws1.Columns.Item(3).ColumnWidth = 5.67
